# Mise à jour de l'application
# Appends new wellness-tracking entries (rows 463-478) for date 08/10/2025 (serial 45938)
# and refreshes the selected cell, matching the authored workbook update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template rows already present on the sheet, used purely so the new rows inherit
# the exact same cell styles (date format, font, centered-empty style, ...).
$blankLocTemplate = $ws.Range("A452:H452")   # G empty -> style s=2
$textLocTemplate  = $ws.Range("A453:H453")   # G has text -> style s=1

# "Adducteur" as used elsewhere in the sheet ends with a non-breaking space (U+00A0).
$nbsp = [char]0x00A0
$adducteurNbsp = "Adducteur" + $nbsp

# New data rows (session of 08/10/2025)
$rows = @(
  @{ r=463; name="Levy Ndoutoume";   c=75; d=7; e=7; f=3; g="Ischio";          h=5  },
  @{ r=464; name="Jeremie Laurent";  c=75; d=8; e=6; f=1; g="Cheville coup";   h=8  },
  @{ r=465; name="Mattheo Haon";     c=75; d=7; e=7; f=0; g="";                h=7  },
  @{ r=466; name="Sofiane Belle";    c=75; d=5; e=6; f=0; g="";                h=4  },
  @{ r=467; name="Ilyes Boughanmi";  c=75; d=6; e=5; f=0; g="";                h=0  },
  @{ r=468; name="Romain Thunet";    c=75; d=5; e=5; f=2; g="Courbature";      h=7  },
  @{ r=469; name="Yoann Martelat";   c=75; d=5; e=6; f=4; g="Genou";           h=6  },
  @{ r=470; name="Kamal Bafounta";   c=75; d=7; e=6; f=2; g="Cheville";        h=7  },
  @{ r=471; name="Maé Clavel";       c=75; d=7; e=7; f=5; g="Quadri";          h=6  },
  @{ r=472; name="Naim Ighbane";     c=75; d=4; e=6; f=0; g="";                h=7  },
  @{ r=473; name="Omar Benyounes";   c=75; d=5; e=3; f=1; g="Quadri";          h=10 },
  @{ r=474; name="Malik Boussaid";   c=75; d=2; e=1; f=0; g="";                h=10 },
  @{ r=475; name="Yoan Zouma";       c=75; d=5; e=8; f=7; g="Cheville droite"; h=5  },
  @{ r=476; name="Naim Dhib";        c=75; d=6; e=7; f=5; g="Partout";         h=6  },
  @{ r=477; name="Emmanuel Valey";   c=75; d=5; e=6; f=7; g=$adducteurNbsp;    h=6  },
  @{ r=478; name="Amine Taiar";      c=75; d=3; e=4; f=4; g="Genou";           h=7  }
)

foreach ($row in $rows) {
  $r = $row.r
  $dstRange = $ws.Range("A" + $r + ":H" + $r)
  if ([string]::IsNullOrEmpty($row.g)) {
    $blankLocTemplate.Copy($dstRange)
  } else {
    $textLocTemplate.Copy($dstRange)
  }

  $ws.Cells.Item($r, 1).Value = 45938
  $ws.Cells.Item($r, 2).Value = $row.name
  $ws.Cells.Item($r, 3).Value = $row.c
  $ws.Cells.Item($r, 4).Value = $row.d
  $ws.Cells.Item($r, 5).Value = $row.e
  $ws.Cells.Item($r, 6).Value = $row.f
  if (-not [string]::IsNullOrEmpty($row.g)) {
    $ws.Cells.Item($r, 7).Value = $row.g
  }
  $ws.Cells.Item($r, 8).Value = $row.h
}

# Extend the Charge (=Volume*Intensité) shared formula down through the new rows.
$ws.Range("I463:I478").Formula = "=C463*D463"

# Match the saved selection state captured in the workbook.
$ws.Range("K473").Select()
